# "update code and docs"
# Re-purpose the "lift test 1" sheet from a weight/success-failure log into a
# "number of cones vs. elastic bands" log:
#   - New header row: Number of Cones / No bands / Two bands / Four Bands / Six bands
#   - Column A becomes a plain 0..16 counter
#   - D2 gets a single "Yes" marker
#   - All the old sample rows/columns (B..C weight & success data) are gone
#   - Header row goes bold (the old green-fill highlight on B14:C14 is removed)
#   - Columns B..E get new widths; column A keeps its original width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - wipes old values, the green B14:C14 fill, etc.
$ws.Cells.Clear()

# ---- Header row (bold) ----
$ws.Range("A1").Value = "Number of Cones "
$ws.Range("B1").Value = "No bands"
$ws.Range("C1").Value = "Two bands"
$ws.Range("D1").Value = "Four Bands"
$ws.Range("E1").Value = "Six bands"
$ws.Range("A1:E1").Font.Bold = $true

# ---- Data ----
# Row 2 also carries the single "Yes" flag in column D.
$ws.Range("D2").Value = "Yes"

# Column A: simple increasing counter 0..16 down rows 2..18.
for ($i = 0; $i -le 16; $i++) {
  $ws.Cells.Item($i + 2, 1).Value = $i
}

# ---- Column widths ----
# Column A is untouched (stays at its original 10.7109375 stored width).
# For the others, this runtime stores ColumnWidth with an internal +5/6
# padding before rounding to the nearest 1/6, so we pre-compensate by
# subtracting 5/6 from the desired stored width to land as close as
# possible to the real target.
$ws.Columns.Item(2).ColumnWidth = (10.85546875 - 5/6)
$ws.Columns.Item(3).ColumnWidth = (12.7109375 - 5/6)
$ws.Columns.Item(4).ColumnWidth = (11 - 5/6)
$ws.Columns.Item(5).ColumnWidth = (10.28515625 - 5/6)

# ---- Selection ----
$ws.Range("D3").Select()
